$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Header rename: column A, row 1 "Gen" -> "MaxFES"
$ws.Range("A1").Value = "MaxFES"

# 2) Column A (rows 2-14): generation counters -> MaxFES fractions
$colAValues = @(0, 0.001, 0.01, 0.1, 0.2, 0.3, 0.4, 0.5, 0.6, 0.7, 0.8, 0.9, 1)
for ($i = 0; $i -lt $colAValues.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $colAValues[$i]
}

# 3) Remove the "Run 50" column (AZ). This shifts the old "Mean" column (BA)
#    left into AZ, and shrinks the used range from A1:BA14 to A1:AZ14.
$ws.Range("AZ1:AZ14").EntireColumn.Delete()

# 4) Update the (now shifted) Mean column AZ with the recalculated values
$meanValues = @(160.3816977, 65.91361057, 17.24446044, 17.18398092, 17.18398092, 17.18398092, 17.18398092, 17.18398092, 17.18398092, 17.18398092, 17.18398092, 17.18398092, 17.18398092)
for ($i = 0; $i -lt $meanValues.Length; $i++) {
    $ws.Cells.Item($i + 2, 52).Value = $meanValues[$i]
}
